$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row 1: new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$headerRange = $ws.Range("I1:J1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

# Data for columns I (I0) and J (IF), rows 2-15
$values = @(
    @(6, 6),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(4, 4),
    @(5, 5),
    @(8, 8),
    @(5, 5),
    @(8, 8),
    @(8, 9),
    @(7, 8),
    @(9, 9),
    @(6, 7),
    @(6, 6)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
